# Auto-generated cell updates derived from the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "60.573.86"
Set-TextValue $ws.Range("E2") "  -2.55%  "
Set-TextValue $ws.Range("D3") "2.903.09"
Set-TextValue $ws.Range("E3") "  -3.69%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.10%  "
Set-TextValue $ws.Range("D5") "586.71"
Set-TextValue $ws.Range("E5") "  -1.33%  "
Set-TextValue $ws.Range("D6") "148.00"
Set-TextValue $ws.Range("E6") "  +0.63%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("E8") "  -2.57%  "
Set-TextValue $ws.Range("D9") "2.905.44"
Set-TextValue $ws.Range("E9") "  -3.58%  "
Set-TextValue $ws.Range("D10") "6.72"
Set-TextValue $ws.Range("E10") "  +6.33%  "
Set-TextValue $ws.Range("E11") "  -3.17%  "
Set-TextValue $ws.Range("E12") "  -2.44%  "
Set-TextValue $ws.Range("D13") "0.0000224"
Set-TextValue $ws.Range("E13") "  -3.17%  "
Set-TextValue $ws.Range("E14") "  -0.33%  "
Set-TextValue $ws.Range("E15") "  +0.41%  "
Set-TextValue $ws.Range("D16") "3.384.98"
Set-TextValue $ws.Range("E16") "  -3.64%  "
Set-TextValue $ws.Range("D17") "6.83"
Set-TextValue $ws.Range("E17") "  -2.03%  "
Set-TextValue $ws.Range("D18") "60.566.17"
Set-TextValue $ws.Range("E18") "  -2.47%  "
Set-TextValue $ws.Range("D19") "2.902.65"
Set-TextValue $ws.Range("E19") "  -3.67%  "
Set-TextValue $ws.Range("D20") "426.39"
Set-TextValue $ws.Range("E21") "  -3.78%  "
Set-TextValue $ws.Range("D22") "0.672"
Set-TextValue $ws.Range("E22") "  -2.41%  "
Set-TextValue $ws.Range("D23") "7.12"
Set-TextValue $ws.Range("E23") "  -3.83%  "
Set-TextValue $ws.Range("D24") "80.61"
Set-TextValue $ws.Range("E24") "  -1.86%  "
Set-TextValue $ws.Range("D25") "11.06"
Set-TextValue $ws.Range("E25") "  +1.60%  "
Set-TextValue $ws.Range("E26") "  -1.52%  "
Set-TextValue $ws.Range("D27") "11.79"
Set-TextValue $ws.Range("E27") "  -2.37%  "
Set-TextValue $ws.Range("E28") "  +0.05%  "
Set-TextValue $ws.Range("D29") "7.28"
Set-TextValue $ws.Range("E29") "  +1.55%  "
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.11%  "
Set-TextValue $ws.Range("E31") "  +3.55%  "
Set-TextValue $ws.Range("E32") "  -3.00%  "
Set-TextValue $ws.Range("D33") "26.52"
Set-TextValue $ws.Range("E33") "  -3.32%  "
Set-TextValue $ws.Range("D34") "0.106"
Set-TextValue $ws.Range("E34") "  -3.29%  "
Set-TextValue $ws.Range("D35") "0.0₃0835"
Set-TextValue $ws.Range("E35") "  -1.93%  "
Set-TextValue $ws.Range("E36") "  -1.58%  "
Set-TextValue $ws.Range("E37") "  -2.64%  "
Set-TextValue $ws.Range("B38") "Stacks"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "2.03"
Set-TextValue $ws.Range("E38") "  -0.53%  "
Set-TextValue $ws.Range("B39") "dogwifhat"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D39") "2.97"
Set-TextValue $ws.Range("E39") "  +0.43%  "
Set-TextValue $ws.Range("D40") "49.36"
Set-TextValue $ws.Range("E40") "  -1.65%  "
Set-TextValue $ws.Range("E41") "  -3.96%  "
Set-TextValue $ws.Range("E42") "  -1.82%  "
Set-TextValue $ws.Range("D43") "0.290"
Set-TextValue $ws.Range("E43") "  +2.28%  "
Set-TextValue $ws.Range("D44") "41.46"
Set-TextValue $ws.Range("E44") "  +1.00%  "
Set-TextValue $ws.Range("E45") "  -1.64%  "
Set-TextValue $ws.Range("D46") "370.14"
Set-TextValue $ws.Range("E46") "  -6.14%  "
Set-TextValue $ws.Range("D47") "133.25"
Set-TextValue $ws.Range("E47") "  -1.06%  "
Set-TextValue $ws.Range("D48") "2.654.12"
Set-TextValue $ws.Range("E49") "  -0.04%  "
Set-TextValue $ws.Range("D50") "25.04"
Set-TextValue $ws.Range("E50") "  +5.59%  "
Set-TextValue $ws.Range("E51") "  -1.09%  "
